# Colocando header nos gráficos
$wb = $excel.ActiveWorkbook

# --- Sheets 1-4: Potencia Acumulada, Geracao Periodo Medio, Atendimento a Ponta, Potencia Incremental ---
$sheetNames14 = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)",
    "Potencia Incremental - SIN(MW)"
)

foreach ($name in $sheetNames14) {
    $ws = $wb.Worksheets.Item($name)

    # Add header in A1, matching the existing header style used by B1:E1
    $ws.Range("B1").Copy()
    $ws.Range("A1").PasteSpecial(-4122)
    $ws.Range("A1").Value = "Fonte/Tecnologia"

    # Fix accented labels in A2:A12 and strip their bold/border styling,
    # matching the plain style already used by the numeric columns
    $ws.Range("A2").Value = "Hidro"
    $ws.Range("A3").Value = "Gás Natural"
    $ws.Range("A4").Value = "Carvão"
    $ws.Range("A5").Value = "Nuclear"
    $ws.Range("A6").Value = "Óleos Comb"
    $ws.Range("A7").Value = "Biomassa"
    $ws.Range("A8").Value = "Eólica"
    $ws.Range("A9").Value = "Solar"
    $ws.Range("A10").Value = "Outros"
    $ws.Range("A11").Value = "Pot. Compl."
    $ws.Range("A12").Value = "GD"

    $ws.Range("B2").Copy()
    $ws.Range("A2:A12").PasteSpecial(-4122)
}

# --- Sheet 5: Emissoes Totais (MtCO2eq) ---
$ws = $wb.Worksheets.Item("Emissoes Totais (MtCO2eq)")

$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Range("A1").Value = "Período"

$ws.Range("A2").Value = "P.Médio"
$ws.Range("A3").Value = "P.Crítico"

$ws.Range("B2").Copy()
$ws.Range("A2:A3").PasteSpecial(-4122)

# Remove row 4 ("Teto") entirely
$ws.Rows.Item(4).Delete()

# --- Sheet 6: Custo Total (bilhões de R$) ---
$ws = $wb.Worksheets.Item("Custo Total (bilhões de R$)")

$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Range("A1").Value = "Tipo Expansão"

$ws.Range("B1").Value = "'2015"

$ws.Range("A2").Value = "Expansão Centralizada"
$ws.Range("B2").Value = 588

$ws.Range("A3").Value = "Expansão por GD"
$ws.Range("B3").Value = 99

$ws.Range("B2").Copy()
$ws.Range("A2:A3").PasteSpecial(-4122)
